# Auto-generated Excel COM-interop script to update Atomos_Profits sheets
# Applies numeric corrections to currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3052
$ws.Cells.Item(98, 9).Value = 3112.8
$ws.Cells.Item(98, 10).Value = 2900
$ws.Cells.Item(98, 11).Value = 3112.8
$ws.Cells.Item(98, 12).Value = 2900
$ws.Cells.Item(98, 13).Value = -1614.8
$ws.Cells.Item(98, 14).Value = -5896
$ws.Cells.Item(99, 8).Value = 1553.6
$ws.Cells.Item(99, 9).Value = 349.2
$ws.Cells.Item(99, 11).Value = 1047.6
$ws.Cells.Item(99, 13).Value = 450.4000000000001
$ws.Cells.Item(111, 8).Value = 1307.1428
$ws.Cells.Item(111, 9).Value = 1223.8
$ws.Cells.Item(111, 10).Value = 1515.5
$ws.Cells.Item(111, 11).Value = 3671.4
$ws.Cells.Item(111, 12).Value = 4546.5
$ws.Cells.Item(111, 13).Value = -604.3999999999996
$ws.Cells.Item(111, 14).Value = -10680.5
$ws.Cells.Item(112, 8).Value = 1613.415
$ws.Cells.Item(112, 10).Value = 1413.6731
$ws.Cells.Item(112, 12).Value = 4241.0193
$ws.Cells.Item(112, 14).Value = -6457.0193
$ws.Cells.Item(116, 8).Value = 3862.8823
$ws.Cells.Item(116, 9).Value = 3199.9
$ws.Cells.Item(116, 11).Value = 3199.9
$ws.Cells.Item(116, 13).Value = 242.0999999999999
$ws.Cells.Item(122, 8).Value = 3052
$ws.Cells.Item(122, 9).Value = 3112.8
$ws.Cells.Item(122, 10).Value = 2900
$ws.Cells.Item(122, 11).Value = 9338.400000000001
$ws.Cells.Item(122, 12).Value = 8700
$ws.Cells.Item(122, 13).Value = -6888.400000000001
$ws.Cells.Item(122, 14).Value = -13600
$ws.Cells.Item(125, 8).Value = 1216
$ws.Cells.Item(125, 9).Value = 1112.5714
$ws.Cells.Item(125, 11).Value = 10013.1426
$ws.Cells.Item(125, 13).Value = -7553.142600000001
$ws.Cells.Item(137, 8).Value = 1854882.6
$ws.Cells.Item(137, 9).Value = 2634709.5
$ws.Cells.Item(137, 11).Value = 7904128.5
$ws.Cells.Item(137, 13).Value = -7901578.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2222.889
$ws.Cells.Item(2, 9).Value = 2168.6667
$ws.Cells.Item(2, 10).Value = 2250
$ws.Cells.Item(2, 11).Value = 2168.6667
$ws.Cells.Item(2, 12).Value = 2250
$ws.Cells.Item(2, 13).Value = -2055.6667
$ws.Cells.Item(2, 14).Value = -2476
$ws.Cells.Item(7, 8).Value = 40000
$ws.Cells.Item(7, 10).Value = 40000
$ws.Cells.Item(7, 12).Value = 40000
$ws.Cells.Item(7, 14).Value = -40228
$ws.Cells.Item(74, 8).Value = 1170.2
$ws.Cells.Item(74, 9).Value = 1199.5294
$ws.Cells.Item(74, 11).Value = 1199.5294
$ws.Cells.Item(74, 13).Value = -325.5293999999999
$ws.Cells.Item(77, 8).Value = 1170.2
$ws.Cells.Item(77, 9).Value = 1199.5294
$ws.Cells.Item(77, 11).Value = 5997.646999999999
$ws.Cells.Item(77, 13).Value = -1629.646999999999
$ws.Cells.Item(116, 8).Value = 2222.889
$ws.Cells.Item(116, 9).Value = 2168.6667
$ws.Cells.Item(116, 10).Value = 2250
$ws.Cells.Item(116, 11).Value = 2168.6667
$ws.Cells.Item(116, 12).Value = 2250
$ws.Cells.Item(116, 13).Value = 125.3332999999998
$ws.Cells.Item(116, 14).Value = -6838
$ws.Cells.Item(132, 8).Value = 1827.2258
$ws.Cells.Item(132, 9).Value = 1490.9584
$ws.Cells.Item(132, 10).Value = 2980.1428
$ws.Cells.Item(132, 11).Value = 4472.8752
$ws.Cells.Item(132, 12).Value = 8940.428400000001
$ws.Cells.Item(132, 13).Value = -1942.8752
$ws.Cells.Item(132, 14).Value = -14000.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2222.889
$ws.Cells.Item(3, 9).Value = 2168.6667
$ws.Cells.Item(3, 10).Value = 2250
$ws.Cells.Item(3, 11).Value = 2168.6667
$ws.Cells.Item(3, 12).Value = 2250
$ws.Cells.Item(3, 13).Value = -2054.6667
$ws.Cells.Item(3, 14).Value = -2478
$ws.Cells.Item(80, 8).Value = 556
$ws.Cells.Item(80, 9).Value = 606.5454999999999
$ws.Cells.Item(80, 11).Value = 606.5454999999999
$ws.Cells.Item(80, 13).Value = 391.4545000000001
$ws.Cells.Item(83, 8).Value = 556
$ws.Cells.Item(83, 9).Value = 606.5454999999999
$ws.Cells.Item(83, 11).Value = 3032.7275
$ws.Cells.Item(83, 13).Value = 1959.2725
$ws.Cells.Item(94, 8).Value = 1200
$ws.Cells.Item(94, 9).Value = 1200
$ws.Cells.Item(94, 11).Value = 1200
$ws.Cells.Item(94, 13).Value = -749
$ws.Cells.Item(134, 8).Value = 2122.9546
$ws.Cells.Item(134, 9).Value = 1627.1818
$ws.Cells.Item(134, 11).Value = 4881.5454
$ws.Cells.Item(134, 13).Value = -2346.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 68.416664
$ws.Cells.Item(7, 9).Value = 54.8
$ws.Cells.Item(7, 10).Value = 78.14286
$ws.Cells.Item(7, 11).Value = 54.8
$ws.Cells.Item(7, 12).Value = 78.14286
$ws.Cells.Item(7, 13).Value = 58.2
$ws.Cells.Item(7, 14).Value = -304.14286
$ws.Cells.Item(107, 8).Value = 1160.56
$ws.Cells.Item(107, 9).Value = 413.5
$ws.Cells.Item(107, 10).Value = 3081.5715
$ws.Cells.Item(107, 11).Value = 413.5
$ws.Cells.Item(107, 12).Value = 3081.5715
$ws.Cells.Item(107, 13).Value = 1506.5
$ws.Cells.Item(107, 14).Value = -6921.5715
$ws.Cells.Item(134, 8).Value = 5714.2856
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 5714.2856
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 17142.8568
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(134, 14).Value = -22212.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 472.33334
$ws.Cells.Item(4, 9).Value = 77
$ws.Cells.Item(4, 10).Value = 1263
$ws.Cells.Item(4, 11).Value = 231
$ws.Cells.Item(4, 12).Value = 3789
$ws.Cells.Item(4, 13).Value = -119
$ws.Cells.Item(4, 14).Value = -4013
$ws.Cells.Item(113, 8).Value = 4348668.5
$ws.Cells.Item(113, 9).Value = 33333626
$ws.Cells.Item(113, 10).Value = 925.05
$ws.Cells.Item(113, 11).Value = 100000878
$ws.Cells.Item(113, 12).Value = 2775.15
$ws.Cells.Item(113, 13).Value = -99998708
$ws.Cells.Item(113, 14).Value = -7115.15
$ws.Cells.Item(122, 8).Value = 1754.909
$ws.Cells.Item(122, 9).Value = 1140.8
$ws.Cells.Item(122, 10).Value = 2266.6667
$ws.Cells.Item(122, 11).Value = 10267.2
$ws.Cells.Item(122, 12).Value = 20400.0003
$ws.Cells.Item(122, 13).Value = -7817.199999999999
$ws.Cells.Item(122, 14).Value = -25300.0003
$ws.Cells.Item(131, 8).Value = 1174.9796
$ws.Cells.Item(131, 9).Value = 11000
$ws.Cells.Item(131, 10).Value = 970.2917
$ws.Cells.Item(131, 11).Value = 33000
$ws.Cells.Item(131, 12).Value = 2910.8751
$ws.Cells.Item(131, 13).Value = -27960
$ws.Cells.Item(131, 14).Value = -12990.8751

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 25063.564
$ws.Cells.Item(102, 9).Value = 2665.4583
$ws.Cells.Item(102, 11).Value = 2665.4583
$ws.Cells.Item(102, 13).Value = -1043.4583
$ws.Cells.Item(107, 8).Value = 889.44446
$ws.Cells.Item(107, 9).Value = 292
$ws.Cells.Item(107, 10).Value = 1828.2858
$ws.Cells.Item(107, 11).Value = 292
$ws.Cells.Item(107, 12).Value = 1828.2858
$ws.Cells.Item(107, 13).Value = 1628
$ws.Cells.Item(107, 14).Value = -5668.2858
$ws.Cells.Item(113, 8).Value = 1806.0454
$ws.Cells.Item(113, 9).Value = 1358.3125
$ws.Cells.Item(113, 10).Value = 3000
$ws.Cells.Item(113, 11).Value = 1358.3125
$ws.Cells.Item(113, 12).Value = 3000
$ws.Cells.Item(113, 13).Value = 811.6875
$ws.Cells.Item(113, 14).Value = -7340
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 1423013.9
$ws.Cells.Item(126, 9).Value = 4133496.5
$ws.Cells.Item(126, 10).Value = 3237.1904
$ws.Cells.Item(126, 11).Value = 12400489.5
$ws.Cells.Item(126, 12).Value = 9711.5712
$ws.Cells.Item(126, 13).Value = -12398019.5
$ws.Cells.Item(126, 14).Value = -14651.5712

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2001252.9
$ws.Cells.Item(7, 9).Value = 2942077
$ws.Cells.Item(7, 10).Value = 2001.5625
$ws.Cells.Item(7, 11).Value = 2942077
$ws.Cells.Item(7, 12).Value = 2001.5625
$ws.Cells.Item(7, 13).Value = -2941965
$ws.Cells.Item(7, 14).Value = -2225.5625
$ws.Cells.Item(18, 8).Value = 80003
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(20, 8).Value = 38003
$ws.Cells.Item(20, 10).Value = 38003
$ws.Cells.Item(20, 12).Value = 38003
$ws.Cells.Item(20, 14).Value = -38455
$ws.Cells.Item(22, 8).Value = 83334530
$ws.Cells.Item(22, 9).Value = 111111496
$ws.Cells.Item(22, 10).Value = 3599.6667
$ws.Cells.Item(22, 11).Value = 111111496
$ws.Cells.Item(22, 12).Value = 3599.6667
$ws.Cells.Item(22, 13).Value = -111111201
$ws.Cells.Item(22, 14).Value = -4189.6667
$ws.Cells.Item(27, 8).Value = 83334530
$ws.Cells.Item(27, 9).Value = 111111496
$ws.Cells.Item(27, 10).Value = 3599.6667
$ws.Cells.Item(27, 11).Value = 111111496
$ws.Cells.Item(27, 12).Value = 3599.6667
$ws.Cells.Item(27, 13).Value = -111111389
$ws.Cells.Item(27, 14).Value = -3813.6667
$ws.Cells.Item(126, 8).Value = 2001252.9
$ws.Cells.Item(126, 9).Value = 2942077
$ws.Cells.Item(126, 10).Value = 2001.5625
$ws.Cells.Item(126, 11).Value = 8826231
$ws.Cells.Item(126, 12).Value = 6004.6875
$ws.Cells.Item(126, 13).Value = -8823761
$ws.Cells.Item(126, 14).Value = -10944.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 18566.666
$ws.Cells.Item(96, 10).Value = 35333.332
$ws.Cells.Item(96, 12).Value = 35333.332
$ws.Cells.Item(96, 14).Value = -38079.332
$ws.Cells.Item(122, 8).Value = 529367.5
$ws.Cells.Item(122, 9).Value = 1113756.9
$ws.Cells.Item(122, 10).Value = 3417
$ws.Cells.Item(122, 11).Value = 3341270.7
$ws.Cells.Item(122, 12).Value = 10251
$ws.Cells.Item(122, 13).Value = -3338820.7
$ws.Cells.Item(122, 14).Value = -15151
$ws.Cells.Item(126, 8).Value = 2858916.8
$ws.Cells.Item(126, 9).Value = 1373.6086
$ws.Cells.Item(126, 11).Value = 4120.825800000001
$ws.Cells.Item(126, 13).Value = -1650.825800000001
$ws.Cells.Item(132, 8).Value = 340204.6
$ws.Cells.Item(132, 9).Value = 529586.3
$ws.Cells.Item(132, 10).Value = 13090.728
$ws.Cells.Item(132, 11).Value = 1588758.9
$ws.Cells.Item(132, 12).Value = 39272.18399999999
$ws.Cells.Item(132, 13).Value = -1586228.9
$ws.Cells.Item(132, 14).Value = -44332.18399999999
